$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cancel-payment report template update: the "amount before tax"
# ("ยอดเงินก่อนภาษี") header cell becomes a "service value"
# ("มูลค่าบริการ") header cell.
$ws.Range("J6").Value = "มูลค่าบริการ"

# Leave the cursor on I9, matching the saved selection in the updated
# template.
$ws.Range("I9").Select() | Out-Null
